# Apply the "new TPM" data update to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 2 (Sending=ECs, Ligand=Vip, Receptor=Adcyap1r1, Target=ECs) ----
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3882076666666667
$ws.Range("H2").Value = 1.164623
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2727436666666667
$ws.Range("N2").Value = 0.8182309999999999
$ws.Range("O2").Value = 0.0516881753217707
$ws.Range("P2").Value = 0.0516881753217707
$ws.Range("Q2").Value = 0.1058811824347778
$ws.Range("R2").Value = 0.9529306419129999
$ws.Range("S2").Value = 0.0516881753217707
$ws.Range("T2").Value = 0.0516881753217707

# ---- Row 3 (Target=FAPs) ----
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3882076666666667
$ws.Range("H3").Value = 1.164623
$ws.Range("O3").Value = 0.7013451998115746
$ws.Range("P3").Value = 0.7013451998115746
$ws.Range("Q3").Value = 1.436677897579556
$ws.Range("R3").Value = 12.930101078216
$ws.Range("S3").Value = 0.7013451998115746
$ws.Range("T3").Value = 0.7013451998115746

# ---- Row 4 (Target cluster becomes "Inflammatory-Mac"; all new data) ----
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3882076666666667
$ws.Range("H4").Value = 1.164623
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02191066666666667
$ws.Range("N4").Value = 0.065732
$ws.Range("O4").Value = 0.004152332458988515
$ws.Range("P4").Value = 0.004152332458988515
$ws.Range("Q4").Value = 0.008505888781777778
$ws.Range("R4").Value = 0.07655299903599999
$ws.Range("S4").Value = 0.004152332458988515
$ws.Range("T4").Value = 0.004152332458988515

# ---- Row 5 (new row; Target=MuSCs, the data that used to live in row 4) ----
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vip"
$ws.Range("C5").Value = "Adcyap1r1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3882076666666667
$ws.Range("H5").Value = 1.164623
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.281261333333333
$ws.Range("N5").Value = 3.843784
$ws.Range("O5").Value = 0.2428142924076661
$ws.Range("P5").Value = 0.2428142924076661
$ws.Range("Q5").Value = 0.4973954726035556
$ws.Range("R5").Value = 4.476559253432
$ws.Range("S5").Value = 0.2428142924076661
$ws.Range("T5").Value = 0.2428142924076661
